$wb = $excel.ActiveWorkbook

# WEB_UI sheet: clear the stray result value that was pasted into D2
$wsWeb = $wb.Worksheets.Item("WEB_UI")
$wsWeb.Range("D2").ClearContents() | Out-Null
$wsWeb.Range("D2").Select() | Out-Null

# WIN_UI sheet: reset E2 back to a blank placeholder value
$wsWin = $wb.Worksheets.Item("WIN_UI")
$wsWin.Range("E2").Value = " "
$wsWin.Range("E2").Select() | Out-Null
